$wb = $excel.ActiveWorkbook

# Sheets that use the "Ano" label pattern for the year header row (B1:E1)
$anoSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Emissoes Totais (MtCO2eq)"
)

foreach ($sheetName in $anoSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("B1").Value = "Ano 2015"
    $ws.Range("C1").Value = "Ano 2030"
    $ws.Range("D1").Value = "Ano 2040"
    $ws.Range("E1").Value = "Ano 2050"
}

# Sheet that uses the "Intervalo" label pattern
$ws4 = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
$ws4.Range("B1").Value = "Intervalo 2015"
$ws4.Range("C1").Value = "Intervalo 2015-2030"
$ws4.Range("D1").Value = "Intervalo 2031-2040"
$ws4.Range("E1").Value = "Intervalo 2041-2050"

# Sheet with only a single year column (B1)
$ws6 = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$ws6.Range("B1").Value = "Ano 2015"
